$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("TC1")
$ws2 = $wb.Worksheets.Item("TC2")

# Set G2 and G6 on both sheets to boolean TRUE
$ws1.Range("G2").Value = $true
$ws1.Range("G6").Value = $true
$ws2.Range("G2").Value = $true
$ws2.Range("G6").Value = $true

# Update the selection shown on TC2 (no longer the selected tab)
$ws2.Activate()
$ws2.Range("G9").Select()

# TC1 becomes the active/selected tab with its own new selection
$ws1.Activate()
$ws1.Range("G8").Select()
